$wb = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Sheet1")
$md  = $wb.Worksheets.Item("MasterData")
$rpt = $wb.Worksheets.Item("Report")

# ---------------------------------------------------------------------------
# 1) Data sheet ("Sheet1" internally / sheet2.xml): add 4 new vehicle-variant
#    rows just above the existing last two rows (which already hold the
#    "BOL PIK UP FB 1.3T - BS6.2 - PS" / "... PS AC" variants).
#    We duplicate those two rows (so formatting + formulas come along for
#    free) and then overwrite two of the four copies with the brand new
#    "Refresh" variants.
# ---------------------------------------------------------------------------
$ws.Rows("35:36").Copy()
$ws.Rows("30:31").Insert()

$ws.Rows("30:31").Copy()
$ws.Rows("32:33").Insert()

# Row 32 -> new variant: BOL PIK UP FB 1.3T - BS6.2 - PS - Refresh
$ws.Range("A32").Value2 = "BOL PIK UP FB 1.3T - BS6.2 - PS - Refresh"
$ws.Range("B32").Value2 = "BOL PIK UP FB 1.3T - BS6.2 - PS - Refresh"
$ws.Range("C32").Value2 = 975001
$ws.Range("E32").Value2 = 52620
$ws.Range("F32").Value2 = 51876
$ws.Range("H32").Value2 = 24375

# Row 33 -> new variant: BOL PIK UP FB 1.3T - BS6.2 - PS AC - Refresh
$ws.Range("A33").Value2 = "BOL PIK UP FB 1.3T - BS6.2 - PS AC - Refresh"
$ws.Range("B33").Value2 = "BOL PIK UP FB 1.3T - BS6.2 - PS AC - Refresh"
$ws.Range("C33").Value2 = 999001
$ws.Range("E33").Value2 = 53405
$ws.Range("F33").Value2 = 53097
$ws.Range("H33").Value2 = 24975

# Re-apply the autofilter / used range over the new extent.
$ws.AutoFilterMode = $false
$ws.Range("A2:T40").AutoFilter()

# _xlnm._FilterDatabase defined name needs to track the new range too.
$fdb = $wb.Names.Item("Sheet1!_FilterDatabase")
$fdb.RefersTo = "=Sheet1!`$A`$2:`$T`$40"

# ---------------------------------------------------------------------------
# 2) MasterData sheet (sheet3.xml): the Vehicle_Model table keeps the same
#    44-row extent, but the variant list is re-ordered so the 4 "BOL PIK UP"
#    variants lead, followed by the 5 "CAMPER" variants, followed by the
#    (duplicate, pre-existing) trailing "BOL PIK UP" PS / PS AC entries.
# ---------------------------------------------------------------------------
$md.Range("A39").Value2 = $md.Range("A35").Value2
$md.Range("A38").Value2 = $md.Range("A34").Value2
$md.Range("A37").Value2 = $md.Range("A33").Value2
$md.Range("A36").Value2 = $md.Range("A32").Value2
$md.Range("A35").Value2 = $md.Range("A31").Value2
$md.Range("A34").Value2 = $md.Range("A30").Value2
$md.Range("A33").Value2 = $md.Range("A29").Value2

$md.Range("A29").Value2 = "BOL PIK UP FB 1.3T - BS6.2 - PS"
$md.Range("A30").Value2 = "BOL PIK UP FB 1.3T - BS6.2 - PS AC"
$md.Range("A31").Value2 = "BOL PIK UP FB 1.3T - BS6.2 - PS - Refresh"
$md.Range("A32").Value2 = "BOL PIK UP FB 1.3T - BS6.2 - PS AC - Refresh"

# ---------------------------------------------------------------------------
# 3) Report sheet: point the variant selector at the new "BOL PIK UP ... PS"
#    variant. This feeds the ANCHORARRAY data-validation list fed by
#    MasterData!C2 (a dynamic array SORT/FILTER/SEARCH/TRANSPOSE formula),
#    so re-assert that formula once D3 / the table data has changed, which
#    makes it recompute and spill across C2:H2 as expected.
# ---------------------------------------------------------------------------
$rpt.Range("D3").Value2 = "BOL PIK UP FB 1.3T - BS6.2 - PS"
$md.Range("C2").Formula = $md.Range("C2").Formula
